$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet: add the new configuration rows (UrlSAC, UsuarioSAC,
# UrlRecConsigConsultas, RutaResultado, RutaDescargas) coming from the
# "Consolidacion de Recaudo" bot configuration.
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")

$wsSettings.Range("A10").Value = "UrlSAC"
$wsSettings.Range("B10").Value = "https://essa-ws12.essa.com.co:9095/GEN/Vistas/Login/Login_Gen.aspx "

$wsSettings.Range("A11").Value = "UsuarioSAC"
$wsSettings.Range("B11").Value = "LARIAGIL"

$wsSettings.Range("A12").Value = "UrlRecConsigConsultas"
$wsSettings.Range("B12").Value = "https://essa-ws12.essa.com.co:9095/SAC/Vistas/App/REC_LOTECO.aspx"

$wsSettings.Range("A13").Value = "RutaResultado"
$wsSettings.Range("A14").Value = "RutaDescargas"
$wsSettings.Range("B14").Value = "C:\Users\jpumarej\Downloads"
$wsSettings.Range("B13").Value = "D:\Leidy\Consolidación recaudos"

# New rows picked up the plain "Normal" formatting used elsewhere on the
# sheet (no wrap, default font) instead of inheriting nothing.
$wsSettings.Range("B10").WrapText = $false
$wsSettings.Range("A11:B12").WrapText = $false

# The header row with the wrapped description text re-autosizes slightly.
$wsSettings.Rows.Item(4).RowHeight = 29

# Column B needed to widen considerably to fit the new URLs.
$wsSettings.Columns.Item(2).EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# Constants sheet: the wrapped description row re-autosizes the same way.
# ---------------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Rows.Item(2).RowHeight = 29

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: Settings becomes the active tab
# (selection parked at B21) and Assets loses the active tab but keeps a
# fresh selection (B15).
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Activate()
$wsAssets.Range("B15").Select()

$wsSettings.Activate()
$wsSettings.Range("B21").Select()
